# Insert a new data row at row 238 (pushes existing rows 238:346 down to 239:347)
# and populate it with a new weekly price observation for "Betarraga".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(238).Insert()

$ws.Range("A238").Value = 4
$ws.Range("B238").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C238").Value = "Los Lagos"
$ws.Range("D238").Value2 = 44806
$ws.Range("E238").Value = 10
$ws.Range("F238").Value = 100114014
$ws.Range("G238").Value = "Betarraga"
$ws.Range("H238").Value = "Sin especificar"
$ws.Range("I238").Value = "Primera"
$ws.Range("J238").Value = 1200
$ws.Range("K238").Value = 1500
$ws.Range("L238").Value = 1500
$ws.Range("M238").Value = 1500
$ws.Range("N238").Value = "$/paquete 5 unidades"
$ws.Range("O238").Value = "Región del Maule"
$ws.Range("P238").Value = 300
$ws.Range("Q238").Value = 5
$ws.Range("R238").Value = "Hortaliza"
